# JIRA:101 - adding drools decision table - add business rules in decision table
#
# Updates the "ACTION" column (E) of the customer-rules decision table with
# real business-rule text (replacing placeholder/test strings), tidies up the
# generated-code template in E8, widens column E / tallens a couple of rows so
# the new, longer text is readable, and moves the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer-rules")
$ws.Activate()

# --- Action column (E) text: replace placeholder strings with real business rules ---
$ws.Range("E10").Value = '"generate a packing slip for shipping"'
$ws.Range("E11").Value = '"create a duplicate packing slip for the royalty department"'
$ws.Range("E12").Value = '"activate that membership"'
$ws.Range("E13").Value = '"apply the upgrade"'
$ws.Range("E14").Value = '"e-mail the owner and inform them of the activation/upgrade"'
$ws.Range("E15").Value = '"add a free ' + [char]0x201C + 'First Aid' + [char]0x201D + ' video to the packing slip"'
$ws.Range("E16").Value = '"generate a commission payment to the agent"'

# Wrap text for the (now much longer) action cells, matching the new column width
$ws.Range("E10:E16").WrapText = $true

# --- Tidy the generated System.out.println snippet in the ACTION header row ---
$ws.Range("E8").Value = "System.out.println(`"Output : `" + `$1 );`norderOutput.setOutput(orderOutput.getOutput()  +`$1);"

# --- Column / row sizing so the new text is readable ---
$ws.Columns("E").ColumnWidth = 163.7142857
$ws.Rows(8).RowHeight = 125.25
$ws.Rows(14).RowHeight = 52.5

# --- Selection / scroll position ---
$ws.Range("E11").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 4
